# =====================================================================
# Edit: add "Player Info" sheet (before ODI Batting), add
# "ODI Batting Extra" sheet (after ODI Batting), and update the
# "ODI Batting" sheet's MATCH_CARD_LINK column into a MATCH_CODE column.
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: apply the workbook's standard bold/centered/bordered header
# style (mirrors the style used for row 1 headers in "ODI Batting").
# ---------------------------------------------------------------------
function Set-HeaderStyle($range) {
    $range.Font.Bold = $true
    $range.HorizontalAlignment = -4108   # xlCenter
    $range.VerticalAlignment = -4160     # xlTop
    $range.Borders.LineStyle = 1
}

# =====================================================================
# 1. "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE
# =====================================================================

$wsOdi = $wb.Worksheets.Item("ODI Batting")

$wsOdi.Range("D1").Value = "MATCH_CODE"

$wsOdi.Range("D2:D26").NumberFormat = "@"
$wsOdi.Range("D2").Value = "3589"
$wsOdi.Range("D3").Value = "3590"
$wsOdi.Range("D4").Value = "3591"
$wsOdi.Range("D5").Value = "3592"
$wsOdi.Range("D6").Value = "3594"
$wsOdi.Range("D7").Value = "3618"
$wsOdi.Range("D8").Value = "3620"
$wsOdi.Range("D9").Value = "3624"
$wsOdi.Range("D10").Value = "3631"
$wsOdi.Range("D11").Value = "3662"
$wsOdi.Range("D12").Value = "3666"
$wsOdi.Range("D13").Value = "3921"
$wsOdi.Range("D14").Value = "3925"
$wsOdi.Range("D15").Value = "3926"
$wsOdi.Range("D16").Value = "3928"
$wsOdi.Range("D17").Value = "3930"
$wsOdi.Range("D18").Value = "3932"
$wsOdi.Range("D19").Value = "3939"
$wsOdi.Range("D20").Value = "3943"
$wsOdi.Range("D21").Value = "3944"
$wsOdi.Range("D22").Value = "3972"
$wsOdi.Range("D23").Value = "3973"
$wsOdi.Range("D24").Value = "3975"
$wsOdi.Range("D25").Value = "3977"
$wsOdi.Range("D26").Value = "3981"

# =====================================================================
# 2. "Player Info" sheet - inserted before "ODI Batting"
# =====================================================================

$wsInfo = $wb.Worksheets.Add($wsOdi)
$wsInfo.Name = "Player Info"

# Re-fetch "ODI Batting" by name: the $wsOdi reference above tracks a
# sheet *position*, and inserting a sheet before it shifts what that
# position now refers to.
$wsOdi = $wb.Worksheets.Item("ODI Batting")

$wsInfo.Range("A1").Value = "ID"
$wsInfo.Range("B1").Value = "NAME"
$wsInfo.Range("C1").Value = "BATTING_HAND"
$wsInfo.Range("D1").Value = "BOWL_STYLE"
Set-HeaderStyle $wsInfo.Range("A1:D1")

$wsInfo.Range("A2").NumberFormat = "@"
$wsInfo.Range("A2").Value = "4252"
$wsInfo.Range("B2").Value = "Sharjeel Khan"
$wsInfo.Range("C2").Value = "Left Handed"
$wsInfo.Range("D2").Value = "Right Arm Leg Break"

# =====================================================================
# 3. "ODI Batting Extra" sheet - inserted after "ODI Batting"
# =====================================================================

$wsExtra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsOdi)
$wsExtra.Name = "ODI Batting Extra"

$wsExtra.Range("A1").Value = "MATCH_CODE"
$wsExtra.Range("B1").Value = "BATTING_POSITION"
$wsExtra.Range("C1").Value = "NUM_4"
$wsExtra.Range("D1").Value = "NUM_6"
$wsExtra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$wsExtra.Range("F1").Value = "MAN_OF_MATCH"
Set-HeaderStyle $wsExtra.Range("A1:F1")

# Pre-format the text-valued columns so numeric-looking strings (match
# codes, counts, percentages) are preserved as literal text, matching
# the source data, instead of being coerced into numbers/percentages.
$wsExtra.Range("A2:A21").NumberFormat = "@"
$wsExtra.Range("C2:C21").NumberFormat = "@"
$wsExtra.Range("D2:D21").NumberFormat = "@"
$wsExtra.Range("E2:E21").NumberFormat = "@"
$wsExtra.Range("F2:F21").NumberFormat = "@"

# Row 2 - MatchCode 3618
$wsExtra.Range("A2").Value = "3618"
$wsExtra.Range("B2").Value = 1
$wsExtra.Range("C2").Value = "4"
$wsExtra.Range("D2").Value = "1"
$wsExtra.Range("E2").Value = "9.15%"
$wsExtra.Range("F2").Value = "NO"

# Row 3 - MatchCode 3620
$wsExtra.Range("A3").Value = "3620"
$wsExtra.Range("F3").Value = "NO"

# Row 4 - MatchCode 3624
$wsExtra.Range("A4").Value = "3624"
$wsExtra.Range("F4").Value = "NO"

# Row 5 - MatchCode 3631
$wsExtra.Range("A5").Value = "3631"
$wsExtra.Range("B5").Value = 1
$wsExtra.Range("C5").Value = "2"
$wsExtra.Range("D5").Value = "0"
$wsExtra.Range("E5").Value = "3.08%"
$wsExtra.Range("F5").Value = "NO"

# Row 6 - MatchCode 3662
$wsExtra.Range("A6").Value = "3662"
$wsExtra.Range("B6").Value = 2
$wsExtra.Range("C6").Value = "1"
$wsExtra.Range("D6").Value = "0"
$wsExtra.Range("E6").Value = "3.86%"
$wsExtra.Range("F6").Value = "NO"

# Row 7 - MatchCode 3666
$wsExtra.Range("A7").Value = "3666"
$wsExtra.Range("B7").Value = 2
$wsExtra.Range("C7").Value = "0"
$wsExtra.Range("D7").Value = "0"
$wsExtra.Range("F7").Value = "NO"

# Row 8 - MatchCode 3921
$wsExtra.Range("A8").Value = "3921"
$wsExtra.Range("B8").Value = 2
$wsExtra.Range("C8").Value = "16"
$wsExtra.Range("D8").Value = "9"
$wsExtra.Range("E8").Value = "45.10%"
$wsExtra.Range("F8").Value = "NO"

# Row 9 - MatchCode 3925
$wsExtra.Range("A9").Value = "3925"
$wsExtra.Range("B9").Value = 2
$wsExtra.Range("C9").Value = "3"
$wsExtra.Range("D9").Value = "0"
$wsExtra.Range("E9").Value = "6.15%"
$wsExtra.Range("F9").Value = "NO"

# Row 10 - MatchCode 3926
$wsExtra.Range("A10").Value = "3926"
$wsExtra.Range("B10").Value = 2
$wsExtra.Range("C10").Value = "0"
$wsExtra.Range("D10").Value = "0"
$wsExtra.Range("F10").Value = "NO"

# Row 11 - MatchCode 3928
$wsExtra.Range("A11").Value = "3928"
$wsExtra.Range("B11").Value = 2
$wsExtra.Range("C11").Value = "12"
$wsExtra.Range("D11").Value = "1"
$wsExtra.Range("E11").Value = "21.09%"
$wsExtra.Range("F11").Value = "NO"

# Row 12 - MatchCode 3930
$wsExtra.Range("A12").Value = "3930"
$wsExtra.Range("F12").Value = "NO"

# Row 13 - MatchCode 3932
$wsExtra.Range("A13").Value = "3932"
$wsExtra.Range("B13").Value = 2
$wsExtra.Range("C13").Value = "2"
$wsExtra.Range("D13").Value = "0"
$wsExtra.Range("E13").Value = "3.29%"
$wsExtra.Range("F13").Value = "NO"

# Row 14 - MatchCode 3939
$wsExtra.Range("A14").Value = "3939"
$wsExtra.Range("B14").Value = 2
$wsExtra.Range("C14").Value = "6"
$wsExtra.Range("D14").Value = "3"
$wsExtra.Range("E14").Value = "19.01%"
$wsExtra.Range("F14").Value = "NO"

# Row 15 - MatchCode 3943
$wsExtra.Range("A15").Value = "3943"
$wsExtra.Range("B15").Value = 2
$wsExtra.Range("C15").Value = "3"
$wsExtra.Range("D15").Value = "1"
$wsExtra.Range("E15").Value = "7.12%"
$wsExtra.Range("F15").Value = "NO"

# Row 16 - MatchCode 3944
$wsExtra.Range("A16").Value = "3944"
$wsExtra.Range("B16").Value = 2
$wsExtra.Range("C16").Value = "5"
$wsExtra.Range("D16").Value = "0"
$wsExtra.Range("E16").Value = "12.34%"
$wsExtra.Range("F16").Value = "NO"

# Row 17 - MatchCode 3972
$wsExtra.Range("A17").Value = "3972"
$wsExtra.Range("B17").Value = 2
$wsExtra.Range("C17").Value = "1"
$wsExtra.Range("D17").Value = "1"
$wsExtra.Range("E17").Value = "10.23%"
$wsExtra.Range("F17").Value = "NO"

# Row 18 - MatchCode 3973
$wsExtra.Range("A18").Value = "3973"
$wsExtra.Range("F18").Value = "NO"

# Row 19 - MatchCode 3975
$wsExtra.Range("A19").Value = "3975"
$wsExtra.Range("B19").Value = 2
$wsExtra.Range("C19").Value = "8"
$wsExtra.Range("D19").Value = "1"
$wsExtra.Range("E19").Value = "19.01%"
$wsExtra.Range("F19").Value = "NO"

# Row 20 - MatchCode 3977
$wsExtra.Range("A20").Value = "3977"
$wsExtra.Range("B20").Value = 2
$wsExtra.Range("C20").Value = "10"
$wsExtra.Range("D20").Value = "3"
$wsExtra.Range("E20").Value = "27.72%"
$wsExtra.Range("F20").Value = "NO"

# Row 21 - MatchCode 3981
$wsExtra.Range("A21").Value = "3981"
$wsExtra.Range("B21").Value = 2
$wsExtra.Range("C21").Value = "9"
$wsExtra.Range("D21").Value = "2"
$wsExtra.Range("E21").Value = "25.32%"
$wsExtra.Range("F21").Value = "NO"
